# CustomerTopListExcel.xlsx - sync with version 26.0.24098.0
# Adds the new ReportMetadata.*/ReportRequest.* helper defined names that
# the report-metadata framework expects, repoints the company-name and
# data-retrieved-date formulas at them, and seeds the sample value used
# by the new ReportRequest.CompanyName lookup.
$wb = $excel.ActiveWorkbook

# 1. Add the new ReportMetadata.* / ReportRequest.* defined names (XLOOKUP
#    helpers against the ReportMetadataValues / ReportRequestValues tables
#    on the hidden "Aggregated Metadata" sheet).
$n = $wb.Names.Add("ReportMetadata.AboutThisReportText", "=_xlfn.XLOOKUP(""About This Report Text"",ReportMetadataValues[[#All],[Report Property]],ReportMetadataValues[[#All],[Report Property Value]],""none"")")
$n.Comment = "Use this function to get the About This Report Text from the ReportMetadataValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportMetadata.AboutThisReportTitle", "=_xlfn.XLOOKUP(""About This Report Title"",ReportMetadataValues[[#All],[Report Property]],ReportMetadataValues[[#All],[Report Property Value]],""none"")")
$n.Comment = "Use this function to get the About This Report Title from the ReportMetadataValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportMetadata.ExtensionName", "=_xlfn.XLOOKUP(""Extension Name"",ReportMetadataValues[[#All],[Report Property]],ReportMetadataValues[[#All],[Report Property Value]],""none"")")
$n.Comment = "Use this function to get the Extension Name from the ReportMetadataValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportMetadata.ExtensionPublisher", "=_xlfn.XLOOKUP(""Extension Publisher"",ReportMetadataValues[[#All],[Report Property]],ReportMetadataValues[[#All],[Report Property Value]],""none"")")
$n.Comment = "Use this function to get the Extension Publisher from the ReportMetadataValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportMetadata.ExtensionVersion", "=_xlfn.XLOOKUP(""Extension Version"",ReportMetadataValues[[#All],[Report Property]],ReportMetadataValues[[#All],[Report Property Value]],""none"")")
$n.Comment = "Use this function to get the Extension Version from the ReportMetadataValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportMetadata.ObjectID", "=_xlfn.XLOOKUP(""Object ID"",ReportMetadataValues[[#All],[Report Property]],ReportMetadataValues[[#All],[Report Property Value]],""none"")")
$n.Comment = "Use this function to get the Object ID from the ReportMetadataValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportMetadata.ObjectName", "=_xlfn.XLOOKUP(""Object Name"",ReportMetadataValues[[#All],[Report Property]],ReportMetadataValues[[#All],[Report Property Value]],""none"")")
$n.Comment = "Use this function to get the Object Name from the ReportMetadataValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportMetadata.ReportHelpLink", "=_xlfn.XLOOKUP(""Report help link"",ReportMetadataValues[[#All],[Report Property]],ReportMetadataValues[[#All],[Report Property Value]],""none"")")
$n.Comment = "Use this function to get the Report help link from the ReportMetadataValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportRequest.CompanyId", "=_xlfn.XLOOKUP(""Company Id"",ReportRequestValues[[#All],[Request Property]],ReportRequestValues[[#All],[Request Property Value]],""none"")")
$n.Comment = "Use this function to get the Company Id from the ReportRequestValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportRequest.CompanyName", "=_xlfn.XLOOKUP(""Company name"",ReportRequestValues[[#All],[Request Property]],ReportRequestValues[[#All],[Request Property Value]],""none"")")
$n.Comment = "Use this function to get the Company name from the ReportRequestValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportRequest.Date", "=_xlfn.XLOOKUP(""Date"",ReportRequestValues[[#All],[Request Property]],ReportRequestValues[[#All],[Request Property Value]],""none"")")
$n.Comment = "Use this function to get the Date from the ReportRequestValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportRequest.EnvironmentName", "=_xlfn.XLOOKUP(""Environment name"",ReportRequestValues[[#All],[Request Property]],ReportRequestValues[[#All],[Request Property Value]],""none"")")
$n.Comment = "Use this function to get the Environment name from the ReportRequestValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportRequest.EnvironmentType", "=_xlfn.XLOOKUP(""Environment type"",ReportRequestValues[[#All],[Request Property]],ReportRequestValues[[#All],[Request Property Value]],""none"")")
$n.Comment = "Use this function to get the Environment type from the ReportRequestValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportRequest.FormatRegion", "=_xlfn.XLOOKUP(""Format Region"",ReportRequestValues[[#All],[Request Property]],ReportRequestValues[[#All],[Request Property Value]],""none"")")
$n.Comment = "Use this function to get the Format Region from the ReportRequestValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportRequest.Language", "=_xlfn.XLOOKUP(""Language"",ReportRequestValues[[#All],[Request Property]],ReportRequestValues[[#All],[Request Property Value]],""none"")")
$n.Comment = "Use this function to get the Language from the ReportRequestValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportRequest.TenantEntraId", "=_xlfn.XLOOKUP(""Tenant Entra Id"",ReportRequestValues[[#All],[Request Property]],ReportRequestValues[[#All],[Request Property Value]],""none"")")
$n.Comment = "Use this function to get the Tenant Entra Id from the ReportRequestValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportRequest.TenantId", "=_xlfn.XLOOKUP(""Tenant Id"",ReportRequestValues[[#All],[Request Property]],ReportRequestValues[[#All],[Request Property Value]],""none"")")
$n.Comment = "Use this function to get the Tenant Id from the ReportRequestValues table in the Aggregated Metadata worksheet"
$n = $wb.Names.Add("ReportRequest.UserName", "=_xlfn.XLOOKUP(""User name"",ReportRequestValues[[#All],[Request Property]],ReportRequestValues[[#All],[Request Property Value]],""none"")")
$n.Comment = "Use this function to get the User name from the ReportRequestValues table in the Aggregated Metadata worksheet"

# 2. Seed the sample/placeholder value for the "Company name" request
#    property so ReportRequest.CompanyName resolves to something sensible.
$wsMeta = $wb.Worksheets.Item("Aggregated Metadata")
$wsMeta.Range("E5").Value = "CompanyNamePlaceholder"

# 3. Point the report-header formulas at the new defined names instead of
#    the inline VLOOKUP calls against ReportRequestValues.
$wsLabel = $wb.Worksheets.Item('$TopCustomerListLabel$')
$wsLabel.Range("M2").Formula = "=ReportRequest.CompanyName"
$wsLabel.Range("M3").Formula = "=IF(ISNA(VLOOKUP(""DataRetrieved"",CaptionData[#All],2,FALSE)),""Data retrieved:"",VLOOKUP(""DataRetrieved"",CaptionData[#All],2,FALSE))&"" ""&TEXT(ReportRequest.Date,""d mmmm yyyy, hh:mm"")"
$wsLabel.Range("M4").Select()
